$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Objetivos:" body text replaced by the teacher info line.
$ws.Range("B10").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C10").Value = "5817650 - Érica Leonor Romão"

# Row 13 previously had no label (A13 blank) and held the teacher info in B/C.
# It now becomes the "Programa resumido:" row with "Semestral" as its value.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows("13").RowHeight = 60

# Row 14 keeps its B/C text but the label becomes "Short syllabus:".
$ws.Range("A14").Value = "Short syllabus:"

# Row 15's label becomes "Programa:" and its value becomes the date string
# (re-using shared string "01/01/2022").
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows("15").RowHeight = 120

# Row 16 keeps its B/C text but the label becomes "Syllabus:".
$ws.Range("A16").Value = "Syllabus:"

# Row 17 loses its B/C text and becomes the lone "Avaliação:" label row.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows("17").RowHeight = 15

# Row 18 gains the "Método:" label plus the reused teacher info text.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C18").Value = "5817650 - Érica Leonor Romão"
$ws.Rows("18").RowHeight = 60

# Row 19 keeps its B/C text but the label becomes "Critério:".
$ws.Range("A19").Value = "Critério:"

# Row 20 keeps its B/C text but the label becomes "Norma de recuperação:".
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21 keeps its B/C text but the label becomes "Bibliografia:".
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows("21").RowHeight = 120

# Row 22 (old "Bibliografia:" + long bibliography text) is removed entirely.
$ws.Rows("22").Delete()
